$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.076.79"
$ws.Range("E2").Value = "  -1.65%  "
$ws.Range("D3").Value = "'1.793.33"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "'313.89"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "'0.5203"
$ws.Range("E7").Value = "  +1.71%  "
$ws.Range("D8").Value = "'0.3815"
$ws.Range("E8").Value = "  -3.64%  "
$ws.Range("D9").Value = "'0.07844"
$ws.Range("E9").Value = "  -4.60%  "
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").Value = "'1.094"
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("D12").Value = "'1.005"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "'6.263"
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").Value = "'20.48"
$ws.Range("E14").Value = "  -3.48%  "
$ws.Range("D15").Value = "'1.792.39"
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").Value = "'7.267"
$ws.Range("E16").Value = "  -3.61%  "
$ws.Range("D17").Value = "'91.97"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "'0.00001083"
$ws.Range("E18").Value = "  -4.23%  "
$ws.Range("D19").Value = "'0.06540"
$ws.Range("E19").Value = "  -1.88%  "
$ws.Range("D20").Value = "'1.004"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").Value = "'17.26"
$ws.Range("E21").Value = "  -3.30%  "
$ws.Range("D22").Value = "'5.933"
$ws.Range("E22").Value = "  -2.83%  "
$ws.Range("D23").Value = "'28.100.31"
$ws.Range("E23").Value = "  -1.68%  "
$ws.Range("D24").Value = "'11.10"
$ws.Range("E24").Value = "  -2.78%  "
$ws.Range("D25").Value = "'2.258"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").Value = "'160.60"
$ws.Range("E26").Value = "  +2.42%  "
$ws.Range("D27").Value = "'20.37"
$ws.Range("E27").Value = "  -4.99%  "
$ws.Range("D28").Value = "'1.991.66"
$ws.Range("E28").Value = "  -2.03%  "
$ws.Range("D29").Value = "'2.320"
$ws.Range("E29").Value = "  -3.73%  "
$ws.Range("D30").Value = "'122.42"
$ws.Range("E30").Value = "  -3.58%  "
$ws.Range("D31").Value = "'0.1068"
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("D32").Value = "'1.047"
$ws.Range("E32").Value = "  -5.95%  "
$ws.Range("D33").Value = "'3.672"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").Value = "'5.527"
$ws.Range("E34").Value = "  -4.37%  "
$ws.Range("D35").Value = "'0.07237"
$ws.Range("E35").Value = "  +2.36%  "
$ws.Range("D36").Value = "'12.21"
$ws.Range("E36").Value = "  +7.88%  "
$ws.Range("D37").Value = "'0.02310"
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("D38").Value = "'8.762"
$ws.Range("E38").Value = "  -1.03%  "
$ws.Range("D39").Value = "'0.2133"
$ws.Range("E39").Value = "  -4.36%  "
$ws.Range("D40").Value = "'5.052"
$ws.Range("E40").Value = "  -4.43%  "
$ws.Range("D41").Value = "'0.6117"
$ws.Range("E41").Value = "  -3.34%  "
$ws.Range("D42").Value = "'1.158"
$ws.Range("E42").Value = "  -2.14%  "
$ws.Range("D43").Value = "'1.372"
$ws.Range("E43").Value = "  -1.87%  "
$ws.Range("D44").Value = "'13.10"
$ws.Range("E44").Value = "  -3.80%  "
$ws.Range("D45").Value = "'3.757"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").Value = "'0.5915"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("D47").Value = "'127.59"
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("D48").Value = "'1.225"
$ws.Range("E48").Value = "  +2.68%  "
$ws.Range("D49").Value = "'1.911"
$ws.Range("E49").Value = "  -4.35%  "
$ws.Range("D50").Value = "'0.06726"
$ws.Range("E50").Value = "  -2.95%  "
$ws.Range("D51").Value = "'72.57"
$ws.Range("E51").Value = "  -2.21%  "
